$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.203.30"
$ws.Range("E2").Value = "  -4.78%  "
$ws.Range("D3").Value = "2.574.98"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.27"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").Value = "2.594.92"
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("E11").Value = "  -5.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.334"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.55%  "
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "3.030.68"
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("D15").Value = "58.187.07"
$ws.Range("E15").Value = "  -4.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.43%  "
$ws.Range("E17").Value = "  -4.46%  "
$ws.Range("D18").Value = "2.588.23"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("E22").Value = "  -4.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.58"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.417"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Value = "2.691.67"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.158"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.38%  "
$ws.Range("D29").Value = "0.0₃0813"
$ws.Range("E29").Value = "  -5.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").ClearFormats()
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.77"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.10"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.960"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("E38").Value = "  -6.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.847"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.04"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "286.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -7.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.57"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("E43").Value = "  -6.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0989"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.610"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.17"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0534"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.81%  "
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0227"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.65"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.96%  "
